$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.165097138958828
$ws.Range("C2").Value = 4.433039901940475
$ws.Range("D2").Value = 0.3262726588111902

$ws.Range("B3").Value = 1.94365635228407
$ws.Range("C3").Value = 2.969999999999999
$ws.Range("D3").Value = 0.5624249451433154

$ws.Range("B4").Value = 1.644922175155494
$ws.Range("C4").Value = 2.115915492957747
$ws.Range("D4").Value = 0.7755311667593128

$ws.Range("B5").Value = 2.167464623991318
$ws.Range("C5").Value = 4.440578147350193
$ws.Range("D5").Value = 0.32332099907314
